$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("time")
$ws.Range("D7").NumberFormat = "m/d/yy"
